$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push the existing "Importance Sampling" block (rows 20-23) and
# everything below it down by 8 rows so the new "Orthogonal Sampling" block
# can be inserted at rows 20-26 (with a blank row 27 before the old block).
$ws.Rows("20:27").Insert()

# The old row 26 (now shifted to row 34) held the huge orthogonal-sampling
# raw-data list string in column A; that data is superseded by the new
# summarized rows below, so drop it.
$ws.Range("A34").ClearContents()

# --- New "Orthogonal Sampling" section ---
$ws.Range("A20").Value = "Orthogonal Sampling"
$ws.Range("A20").Font.Bold = $true

$ws.Range("A21").Value = 10
$ws.Range("B21").Value = 1000
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 1.32
$ws.Range("E21").Value = 1.0047885349664301
$ws.Range("G21").Value = "[1.32, 1.536, 1.5155999999999998, 1.50836, 1.5121200000000001, 1.5105792]"

$ws.Range("A22").Value = 100
$ws.Range("B22").Value = 1000
$ws.Range("C22").Value = 40
$ws.Range("D22").Value = 1.536
$ws.Range("E22").Value = 0.25996922894835001
$ws.Range("G22").Value = "[1.0047885349664378, 0.25996922894835073, 0.11559861590866907, 0.035726942214524875, 0.006726829862572712, 0.003651293600903659]"

$ws.Range("A23").Value = 1000
$ws.Range("B23").Value = 1000
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = 1.5155999999999901
$ws.Range("E23").Value = 0.115598615908669

$ws.Range("A24").Value = 10000
$ws.Range("B24").Value = 1000
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 1.5083599999999999
$ws.Range("E24").Value = 0.035726942214524798

$ws.Range("A25").Value = 100000
$ws.Range("B25").Value = 1000
$ws.Range("C25").Value = 40
$ws.Range("D25").Value = 1.5121199999999999
$ws.Range("E25").Value = 0.0067268298625727096

$ws.Range("A26").Value = 1000000
$ws.Range("B26").Value = 1000
$ws.Range("C26").Value = 40
$ws.Range("D26").Value = 1.5105792
$ws.Range("E26").Value = 0.0036512936009036501

# Row 27 stays blank, matching the spacing used before every other section.

$ws.Range("A35").Select()
